$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.233.59'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '1.589.15'
$ws.Range("E3").Value = '  +0.90%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''212.33'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("D9").Value = '''0.0608'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '''19.37'
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("D11").Value = '''0.0850'
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("D12").Value = '1.812.08'
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("D13").Value = '1.635.65'
$ws.Range("E13").Value = '  +3.79%  '
$ws.Range("D14").Value = '''4.04'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("E15").Value = '  +1.73%  '
$ws.Range("D16").Value = '''64.38'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '26.228.82'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("D19").Value = '''7.41'
$ws.Range("E19").Value = '  +2.06%  '
$ws.Range("D20").Value = '''213.34'
$ws.Range("E20").Value = '  +3.20%  '
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("D23").Value = '''9.00'
$ws.Range("E23").Value = '  +1.69%  '
$ws.Range("D24").Value = '''2.14'
$ws.Range("E24").Value = '  -2.47%  '
$ws.Range("D25").Value = '''143.96'
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").Value = '''7.06'
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("D29").Value = '''15.20'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  -1.61%  '
$ws.Range("D31").Value = '''1.16'
$ws.Range("E31").Value = '  +1.08%  '
$ws.Range("D32").Value = '''3.20'
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").Value = '1.334.58'
$ws.Range("E34").Value = '  +4.36%  '
$ws.Range("E35").Value = '  -0.74%  '
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("E37").Value = '  -3.47%  '
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("E40").Value = '  +3.02%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("E42").Value = '  -7.47%  '
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").Value = '''0.766'
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("D45").Value = '''61.83'
$ws.Range("E45").Value = '  -0.97%  '
$ws.Range("D46").Value = '1.724.03'
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").Value = '''85.47'
$ws.Range("E47").Value = '  -4.07%  '
$ws.Range("E48").Value = '  -2.25%  '
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("E50").Value = '  -2.97%  '
$ws.Range("E51").Value = '  -0.35%  '
